$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price values so Excel does not
# auto-convert them to numbers (the source data stores these as text).
$textForceCells = @('D5', 'D6', 'D9', 'D11', 'D12', 'D17', 'D19', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D28', 'D29', 'D31', 'D32', 'D33', 'D34', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D43', 'D46', 'D47', 'D48', 'D50', 'D51')
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated cell values from the crypto price refresh.
$ws.Range('D2').Value = '42.756.13'
$ws.Range('E2').Value = '  -6.74%  '
$ws.Range('D3').Value = '2.548.69'
$ws.Range('E3').Value = '  -4.62%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '299.87'
$ws.Range('E5').Value = '  -3.81%  '
$ws.Range('D6').Value = '92.54'
$ws.Range('E6').Value = '  -6.03%  '
$ws.Range('E7').Value = '  -4.06%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.548'
$ws.Range('E9').Value = '  -5.55%  '
$ws.Range('E10').Value = '  -6.10%  '
$ws.Range('D11').Value = '0.0807'
$ws.Range('E11').Value = '  -4.94%  '
$ws.Range('D12').Value = '7.72'
$ws.Range('E12').Value = '  -4.99%  '
$ws.Range('E13').Value = '  +5.17%  '
$ws.Range('D14').Value = '2.935.48'
$ws.Range('E14').Value = '  -4.27%  '
$ws.Range('D15').Value = '2.519.52'
$ws.Range('E15').Value = '  -5.24%  '
$ws.Range('E16').Value = '  -5.72%  '
$ws.Range('D17').Value = '14.22'
$ws.Range('E17').Value = '  -6.04%  '
$ws.Range('D18').Value = '42.756.81'
$ws.Range('E18').Value = '  -6.86%  '
$ws.Range('D19').Value = '12.90'
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('D20').Value = '0.0₃0984'
$ws.Range('E20').Value = '  -3.48%  '
$ws.Range('D21').Value = '6.58'
$ws.Range('E21').Value = '  -3.73%  '
$ws.Range('D22').Value = '71.80'
$ws.Range('E22').Value = '  -3.27%  '
$ws.Range('D23').Value = '255.36'
$ws.Range('E23').Value = '  -9.91%  '
$ws.Range('D24').Value = '2.93'
$ws.Range('E24').Value = '  -4.73%  '
$ws.Range('D25').Value = '2.13'
$ws.Range('E25').Value = '  -5.39%  '
$ws.Range('D26').Value = '29.25'
$ws.Range('E26').Value = '  -5.30%  '
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').Value = '10.12'
$ws.Range('E28').Value = '  -3.99%  '
$ws.Range('D29').Value = '37.03'
$ws.Range('E29').Value = '  -3.70%  '
$ws.Range('E30').Value = '  -2.75%  '
$ws.Range('D31').Value = '6.02'
$ws.Range('E31').Value = '  -3.19%  '
$ws.Range('D32').Value = '152.50'
$ws.Range('E32').Value = '  -1.93%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').Value = '2.17'
$ws.Range('E33').Value = '  -8.60%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = '2.76'
$ws.Range('E34').Value = '  -2.37%  '
$ws.Range('E35').Value = '  -10.38%  '
$ws.Range('D36').Value = '0.0795'
$ws.Range('E36').Value = '  -5.61%  '
$ws.Range('D37').Value = '0.115'
$ws.Range('E37').Value = '  -5.18%  '
$ws.Range('D38').Value = '17.16'
$ws.Range('E38').Value = '  +5.98%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').Value = '24.15'
$ws.Range('E39').Value = '  -6.16%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = '0.119'
$ws.Range('E40').Value = '  -4.41%  '
$ws.Range('D41').Value = '3.43'
$ws.Range('E41').Value = '  -4.76%  '
$ws.Range('D42').Value = '0.0311'
$ws.Range('E42').Value = '  -5.35%  '
$ws.Range('D43').Value = '3.88'
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('D44').Value = '2.085.40'
$ws.Range('E44').Value = '  -3.33%  '
$ws.Range('E45').Value = '  -0.22%  '
$ws.Range('D46').Value = '1.65'
$ws.Range('E46').Value = '  +3.50%  '
$ws.Range('D47').Value = '9.07'
$ws.Range('E47').Value = '  -2.87%  '
$ws.Range('D48').Value = '84.53'
$ws.Range('E48').Value = '  -10.51%  '
$ws.Range('D49').Value = '2.790.47'
$ws.Range('E49').Value = '  -4.36%  '
$ws.Range('D50').Value = '105.00'
$ws.Range('E50').Value = '  -6.04%  '
$ws.Range('D51').Value = '1.67'
$ws.Range('E51').Value = '  -4.06%  '
